$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -4
    3  = -5
    4  = 2
    5  = 1
    6  = -2
    7  = 10
    8  = 2
    9  = 1
    10 = -1
    11 = -2
    12 = -4
    13 = 3
    14 = -5
    15 = 5
    16 = 0
    17 = 8
    18 = -3
    19 = -2
    21 = -3
    22 = -1
    23 = -1
    24 = 6
    25 = -3
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
